$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '86.492.54'
$ws.Range('E2').Value = '  +3.10%  '
$ws.Range('D3').Value = '3.267.37'
$ws.Range('E3').Value = '  +1.14%  '
$ws.Range('E4').Value = '  +0.12%  '
$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '211.26'
$c.Style = 'Normal'
$ws.Range('E5').Value = '  -3.19%  '
$c = $ws.Range('D6')
$c.NumberFormat = '@'
$c.Value = '626.94'
$c.Style = 'Normal'
$ws.Range('E6').Value = '  +0.32%  '
$c = $ws.Range('D7')
$c.NumberFormat = '@'
$c.Value = '0.370'
$c.Style = 'Normal'
$ws.Range('E7').Value = '  +20.06%  '
$c = $ws.Range('D8')
$c.NumberFormat = '@'
$c.Value = '0.691'
$c.Style = 'Normal'
$ws.Range('E8').Value = '  +17.37%  '
$ws.Range('E9').Value = '  +0.12%  '
$ws.Range('D10').Value = '3.268.50'
$ws.Range('E10').Value = '  +1.30%  '
$c = $ws.Range('D11')
$c.NumberFormat = '@'
$c.Value = '0.572'
$c.Style = 'Normal'
$ws.Range('E11').Value = '  -4.63%  '
$ws.Range('E12').Value = '  +8.23%  '
$ws.Range('E13').Value = '  -7.15%  '
$ws.Range('E14').Value = '  +4.71%  '
$ws.Range('D15').Value = '3.871.92'
$ws.Range('E15').Value = '  +1.40%  '
$ws.Range('E16').Value = '  -1.13%  '
$ws.Range('D17').Value = '86.370.19'
$ws.Range('E17').Value = '  +3.46%  '
$ws.Range('D18').Value = '3.272.92'
$ws.Range('E18').Value = '  +1.70%  '
$c = $ws.Range('D19')
$c.NumberFormat = '@'
$c.Value = '14.01'
$c.Style = 'Normal'
$ws.Range('E19').Value = '  -2.14%  '
$ws.Range('E20').Value = '  -5.50%  '
$c = $ws.Range('D21')
$c.NumberFormat = '@'
$c.Value = '431.35'
$c.Style = 'Normal'
$ws.Range('E21').Value = '  -3.34%  '
$c = $ws.Range('D22')
$c.NumberFormat = '@'
$c.Value = '8.86'
$c.Style = 'Normal'
$ws.Range('E22').Value = '  -1.91%  '
$ws.Range('E23').Value = '  +2.39%  '
$c = $ws.Range('D24')
$c.NumberFormat = '@'
$c.Value = '7.26'
$c.Style = 'Normal'
$ws.Range('E24').Value = '  -2.21%  '
$c = $ws.Range('D25')
$c.NumberFormat = '@'
$c.Value = '12.52'
$c.Style = 'Normal'
$ws.Range('E25').Value = '  +5.39%  '
$c = $ws.Range('D26')
$c.NumberFormat = '@'
$c.Value = '5.09'
$c.Style = 'Normal'
$ws.Range('E26').Value = '  -1.42%  '
$ws.Range('D27').Value = '3.442.24'
$ws.Range('E27').Value = '  +1.41%  '
$c = $ws.Range('D28')
$c.NumberFormat = '@'
$c.Value = '76.07'
$c.Style = 'Normal'
$ws.Range('E28').Value = '  -2.93%  '
$ws.Range('E29').Value = '  +4.25%  '
$ws.Range('E30').Value = '  +0.27%  '
$ws.Range('E31').Value = '  +12.13%  '
$ws.Range('E32').Value = '  +0.11%  '
$c = $ws.Range('D33')
$c.NumberFormat = '@'
$c.Value = '8.79'
$c.Style = 'Normal'
$ws.Range('E33').Value = '  -3.45%  '
$c = $ws.Range('D34')
$c.NumberFormat = '@'
$c.Value = '541.64'
$c.Style = 'Normal'
$ws.Range('E34').Value = '  -4.61%  '
$c = $ws.Range('D35')
$c.NumberFormat = '@'
$c.Value = '1.42'
$c.Style = 'Normal'
$ws.Range('E35').Value = '  -4.37%  '
$ws.Range('E36').Value = '  -2.61%  '
$c = $ws.Range('D37')
$c.NumberFormat = '@'
$c.Value = '6.96'
$c.Style = 'Normal'
$ws.Range('E37').Value = '  +12.40%  '
$c = $ws.Range('D38')
$c.NumberFormat = '@'
$c.Value = '0.137'
$c.Style = 'Normal'
$ws.Range('E38').Value = '  -10.74%  '
$c = $ws.Range('D39')
$c.NumberFormat = '@'
$c.Value = '22.48'
$c.Style = 'Normal'
$ws.Range('E39').Value = '  -2.69%  '
$c = $ws.Range('D40')
$c.NumberFormat = '@'
$c.Value = '0.999'
$c.Style = 'Normal'
$ws.Range('E40').Value = '  +0.16%  '
$c = $ws.Range('D41')
$c.NumberFormat = '@'
$c.Value = '21.55'
$c.Style = 'Normal'
$ws.Range('E41').Value = '  +3.14%  '
$c = $ws.Range('D42')
$c.NumberFormat = '@'
$c.Value = '0.392'
$c.Style = 'Normal'
$ws.Range('E42').Value = '  -3.65%  '
$ws.Range('E43').Value = '  -2.40%  '
$ws.Range('E44').Value = '  -3.87%  '
$ws.Range('E45').Value = '  -0.01%  '
$c = $ws.Range('D46')
$c.NumberFormat = '@'
$c.Value = '157.35'
$c.Style = 'Normal'
$ws.Range('E46').Value = '  -1.87%  '
$c = $ws.Range('D47')
$c.NumberFormat = '@'
$c.Value = '178.76'
$c.Style = 'Normal'
$ws.Range('E47').Value = '  -4.90%  '
$c = $ws.Range('D48')
$c.NumberFormat = '@'
$c.Value = '44.38'
$c.Style = 'Normal'
$ws.Range('E48').Value = '  -0.85%  '
$ws.Range('E49').Value = '  -1.66%  '
$c = $ws.Range('D50')
$c.NumberFormat = '@'
$c.Value = '4.23'
$c.Style = 'Normal'
$ws.Range('E50').Value = '  +0.04%  '
$ws.Range('B51').Value = 'ARBITRUM'
$ws.Range('C51').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$c = $ws.Range('D51')
$c.NumberFormat = '@'
$c.Value = '0.626'
$c.Style = 'Normal'
$ws.Range('E51').Value = '  -1.37%  '
